# Updates the cryptos worksheet with freshly scraped coinranking.com data.
# Commit: Updated cryptos list on Sun Aug 11 14:28:52 UTC 2024 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    # Force the written value to stay a text cell (matches the source data,
    # which stores prices/percentages/coin names as inline strings) even when
    # the text looks numeric (e.g. "524.92"), then restore the original style
    # so no stray number-format / style is left behind on the cell.
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = $origStyle
}

$data = @(
    ,@(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "60.521.30", "  +0.14%  ")
    ,@(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "2.636.63", "  +1.47%  ")
    ,@(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  -0.06%  ")
    ,@(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "524.92", "  +2.21%  ")
    ,@(6, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "151.26", "  -1.31%  ")
    ,@(7, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "0.998", "  -0.02%  ")
    ,@(8, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.576", "  -3.88%  ")
    ,@(9, "LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "2.646.71", "  +1.38%  ")
    ,@(10, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "6.41", "  -3.34%  ")
    ,@(11, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.106", "  +2.60%  ")
    ,@(12, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.345", "  +0.10%  ")
    ,@(13, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.128", "  -0.65%  ")
    ,@(14, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "3.097.07", "  +1.39%  ")
    ,@(15, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "60.506.16", "  +0.02%  ")
    ,@(16, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "21.60", "  +0.14%  ")
    ,@(17, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000139", "  -0.34%  ")
    ,@(18, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "2.631.88", "  +0.89%  ")
    ,@(19, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "4.67", "  -1.57%  ")
    ,@(20, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "347.69", "  -3.06%  ")
    ,@(21, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "10.51", "  -0.33%  ")
    ,@(22, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "6.20", "  +0.21%  ")
    ,@(23, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.995", "  -0.43%  ")
    ,@(24, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "61.07", "  +0.00%  ")
    ,@(25, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.423", "  -0.70%  ")
    ,@(26, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.164", "  -0.85%  ")
    ,@(27, "Binance-PegBSC-USD", "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd", "1.00", "  +0.23%  ")
    ,@(28, "PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0₃0833", "  +0.02%  ")
    ,@(29, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "7.15", "  -1.23%  ")
    ,@(30, "USDe", "https://coinranking.com/coin/exbfr2U-0+usde-usde", "0.999", "  -0.03%  ")
    ,@(31, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "6.05", "  +2.08%  ")
    ,@(32, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.60", "  +1.38%  ")
    ,@(33, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "19.11", "  -1.45%  ")
    ,@(34, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "150.49", "  +0.10%  ")
    ,@(35, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "4.03", "  +0.45%  ")
    ,@(36, "SuiNetwork", "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui", "0.899", "  -1.98%  ")
    ,@(37, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "1.16", "  -1.97%  ")
    ,@(38, "Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "0.881", "  +4.71%  ")
    ,@(39, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "36.66", "  +1.01%  ")
    ,@(40, "Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "1.46", "  -1.22%  ")
    ,@(41, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "3.68", "  -1.63%  ")
    ,@(42, "Bittensor", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao", "292.67", "  +1.70%  ")
    ,@(43, "Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.634", "  +2.68%  ")
    ,@(44, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.101", "  -0.67%  ")
    ,@(45, "FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "0.998", "  +0.07%  ")
    ,@(46, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "19.99", "  +2.15%  ")
    ,@(47, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.0553", "  -0.03%  ")
    ,@(48, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "4.85", "  -1.89%  ")
    ,@(49, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0236", "  +0.57%  ")
    ,@(50, "WhiteBITCoin", "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt", "10.40", "  +0.94%  ")
    ,@(51, "InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "18.81", "  -1.97%  ")
)

foreach ($row in $data) {
    $r = $row[0]
    Set-TextValue $ws.Cells.Item($r, 2) $row[1]
    Set-TextValue $ws.Cells.Item($r, 3) $row[2]
    Set-TextValue $ws.Cells.Item($r, 4) $row[3]
    Set-TextValue $ws.Cells.Item($r, 5) $row[4]
}

"Updated $($data.Count) rows."
